# germanized assignments, minor typo fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to its German title.
$ws.Name = "Buchverkäufe (Fiktion)"

# Re-format the yearly sales figures using the German/Euro accounting
# number format (what you get from the "Euro" locale accounting format)
# instead of the old US-dollar custom format.
$ws.Range("B2:F6").NumberFormat = '_-* #,##0.00\ [$€-407]_-;\-* #,##0.00\ [$€-407]_-;_-* "-"??\ [$€-407]_-;_-@_-'

# Widen the data columns very slightly to fit the new format, and select
# the re-formatted range (matching the state Excel leaves behind after the
# reformat operation).
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

$ws.Range("B2:F6").Select()
